# Generate Report for Handback
# Populates the "Latest Target File" (F) and "Latest Handback File" (G)
# columns on the zh-cn / de-de report sheets, flips the Status text from
# "Ready for handoff" to "Handed back: in sync with en-US", and stamps the
# handback datetime for each locale.

$wb = $excel.ActiveWorkbook

# Hyperlink colour used throughout this workbook for "handoff/handback"
# style cells (matches the existing custom Hyperlink font: underline + FF6495ED).
$linkColor = 15570276

function Set-HandbackRow($ws, $row, $mdUrl, $xlfUrl, $xlfName, $handbackDate) {

    # Latest Target File -> same source file as the handoff (a.md / b.md)
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = "a.md"
    $ws.Hyperlinks.Add($fCell, $mdUrl, "", "", "a.md")
    $fCell.Font.Underline = $true
    $fCell.Font.Color = $linkColor
    $fCell.Font.Name = "Calibri"

    # Latest Handback File -> the translated .xlf that came back
    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $xlfName
    $ws.Hyperlinks.Add($gCell, $xlfUrl, "", "", $xlfName)
    $gCell.Font.Underline = $true
    $gCell.Font.Color = $linkColor
    $gCell.Font.Name = "Calibri"

    # Latest Handback DateTime
    $ws.Cells.Item($row, 8).Value = $handbackDate

    # Status
    $ws.Cells.Item($row, 3).Value = "Handed back: in sync with en-US"
}

$zhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/362f11a1b20b7ad7e852bde048f7d4c14d519aba/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a7a0fa4ac1c6cf74d96ef9e3a9e4281fec136f6a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$aMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1f829807538153f3cc41cddef8d455babb2a039a/e2e/a.md"

$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZh 2 $aMdUrl $zhUrl $zhXlfName "2016-03-23 08:36:50"
Set-HandbackRow $wsZh 3 $aMdUrl $zhUrl $zhXlfName "2016-03-23 08:36:50"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDe 2 $aMdUrl $deUrl $deXlfName "2016-03-23 08:36:59"
Set-HandbackRow $wsDe 3 $aMdUrl $deUrl $deXlfName "2016-03-23 08:36:59"

$wb.Save()
